$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2679.6
$ws.Range("J17").Value = 2880
$ws.Range("L17").Value = 8640
$ws.Range("N17").Value = -8976
$ws.Range("H28").Value = 785.3333
$ws.Range("I28").Value = 758.5
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 758.5
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -273.5
$ws.Range("N28").Value = -1970
$ws.Range("H32").Value = 6108.533
$ws.Range("I32").Value = 4147.5557
$ws.Range("J32").Value = 9050
$ws.Range("K32").Value = 4147.5557
$ws.Range("L32").Value = 9050
$ws.Range("M32").Value = -3821.5557
$ws.Range("N32").Value = -9702
$ws.Range("H43").Value = 3486.6155
$ws.Range("I43").Value = 1749
$ws.Range("K43").Value = 1749
$ws.Range("M43").Value = -1680
$ws.Range("H92").Value = 667.1429000000001
$ws.Range("I92").Value = 781.64
$ws.Range("J92").Value = 380.9
$ws.Range("K92").Value = 781.64
$ws.Range("L92").Value = 380.9
$ws.Range("M92").Value = 466.36
$ws.Range("N92").Value = -2876.9
$ws.Range("H121").Value = 1150.8
$ws.Range("J121").Value = 1150.8
$ws.Range("L121").Value = 3452.4
$ws.Range("N121").Value = -6946.4
$ws.Range("H133").Value = 89749
$ws.Range("J133").Value = 89749
$ws.Range("L133").Value = 89749
$ws.Range("N133").Value = -99869
$ws.Range("H138").Value = 3022.0864
$ws.Range("J138").Value = 3414.0747
$ws.Range("L138").Value = 10242.2241
$ws.Range("N138").Value = -20522.2241

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2289.1924
$ws.Range("I2").Value = 1880.6666
$ws.Range("K2").Value = 1880.6666
$ws.Range("M2").Value = -1767.6666
$ws.Range("H32").Value = 5376415
$ws.Range("I32").Value = 6258521.5
$ws.Range("J32").Value = 671848
$ws.Range("K32").Value = 6258521.5
$ws.Range("L32").Value = 671848
$ws.Range("M32").Value = -6258234.5
$ws.Range("N32").Value = -672422
$ws.Range("H61").Value = 2541.2
$ws.Range("I61").Value = 2393.25
$ws.Range("K61").Value = 2393.25
$ws.Range("M61").Value = -2181.25
$ws.Range("H88").Value = 1883.5
$ws.Range("I88").Value = 1338.5
$ws.Range("J88").Value = 2519.3333
$ws.Range("K88").Value = 1338.5
$ws.Range("L88").Value = 2519.3333
$ws.Range("M88").Value = -932.5
$ws.Range("N88").Value = -3331.3333
$ws.Range("H91").Value = 1883.5
$ws.Range("I91").Value = 1338.5
$ws.Range("J91").Value = 2519.3333
$ws.Range("K91").Value = 1338.5
$ws.Range("L91").Value = 2519.3333
$ws.Range("M91").Value = 65.5
$ws.Range("N91").Value = -5327.3333
$ws.Range("H116").Value = 2289.1924
$ws.Range("I116").Value = 1880.6666
$ws.Range("K116").Value = 1880.6666
$ws.Range("M116").Value = 413.3334
$ws.Range("H136").Value = 2541.2
$ws.Range("I136").Value = 2393.25
$ws.Range("K136").Value = 7179.75
$ws.Range("M136").Value = -4629.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2289.1924
$ws.Range("I3").Value = 1880.6666
$ws.Range("K3").Value = 1880.6666
$ws.Range("M3").Value = -1766.6666
$ws.Range("H20").Value = 9085.166999999999
$ws.Range("I20").Value = 10517.2
$ws.Range("J20").Value = 1925
$ws.Range("K20").Value = 10517.2
$ws.Range("L20").Value = 1925
$ws.Range("M20").Value = -10270.2
$ws.Range("N20").Value = -2419
$ws.Range("H134").Value = 15768133
$ws.Range("I134").Value = 7578170
$ws.Range("J134").Value = 83335330
$ws.Range("K134").Value = 22734510
$ws.Range("L134").Value = 250005990
$ws.Range("M134").Value = -22731975
$ws.Range("N134").Value = -250011060

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1577.3636
$ws.Range("I16").Value = 1301
$ws.Range("J16").Value = 2314.3333
$ws.Range("K16").Value = 1301
$ws.Range("L16").Value = 2314.3333
$ws.Range("M16").Value = -1014
$ws.Range("N16").Value = -2888.3333
$ws.Range("H22").Value = 141.8
$ws.Range("I22").Value = 141.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 141.8
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 208.2
$ws.Range("H31").Value = 1384.3334
$ws.Range("I31").Value = 1306.5625
$ws.Range("J31").Value = 2006.5
$ws.Range("K31").Value = 1306.5625
$ws.Range("L31").Value = 2006.5
$ws.Range("M31").Value = -1011.5625
$ws.Range("N31").Value = -2596.5
$ws.Range("H34").Value = 1384.3334
$ws.Range("I34").Value = 1306.5625
$ws.Range("J34").Value = 2006.5
$ws.Range("K34").Value = 1306.5625
$ws.Range("L34").Value = 2006.5
$ws.Range("M34").Value = -1104.5625
$ws.Range("N34").Value = -2410.5
$ws.Range("H113").Value = 1577.3636
$ws.Range("I113").Value = 1301
$ws.Range("J113").Value = 2314.3333
$ws.Range("K113").Value = 1301
$ws.Range("L113").Value = 2314.3333
$ws.Range("M113").Value = 869
$ws.Range("N113").Value = -6654.3333
$ws.Range("H122").Value = 6428.909
$ws.Range("I122").Value = 6428.909
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 19286.727
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -16836.727
$ws.Range("H134").Value = 12501830
$ws.Range("I134").Value = 2173.1667
$ws.Range("K134").Value = 6519.500100000001
$ws.Range("M134").Value = -3984.500100000001
$ws.Range("N22").ClearContents()
$ws.Range("N122").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.90000000000001
$ws.Range("I2").Value = 48.333332
$ws.Range("J2").Value = 84.85714
$ws.Range("K2").Value = 289.999992
$ws.Range("L2").Value = 509.14284
$ws.Range("M2").Value = -176.999992
$ws.Range("N2").Value = -735.14284
$ws.Range("H40").Value = 103.333336
$ws.Range("I40").Value = 120
$ws.Range("J40").Value = 20
$ws.Range("K40").Value = 480
$ws.Range("L40").Value = 80
$ws.Range("M40").Value = -411
$ws.Range("N40").Value = -218
$ws.Range("H76").Value = 7500
$ws.Range("I76").Value = 7500
$ws.Range("K76").Value = 22500
$ws.Range("M76").Value = -22117
$ws.Range("H79").Value = 7500
$ws.Range("I79").Value = 7500
$ws.Range("K79").Value = 22500
$ws.Range("M79").Value = -21174
$ws.Range("H113").Value = 529.381
$ws.Range("I113").Value = 202.66667
$ws.Range("J113").Value = 660.06665
$ws.Range("K113").Value = 608.00001
$ws.Range("L113").Value = 1980.19995
$ws.Range("M113").Value = 1561.99999
$ws.Range("N113").Value = -6320.19995
$ws.Range("H129").Value = 2433.8667
$ws.Range("J129").Value = 2707.3076
$ws.Range("L129").Value = 8121.9228
$ws.Range("N129").Value = -18121.9228
$ws.Range("H132").Value = 4028.375
$ws.Range("I132").Value = 1643.3889
$ws.Range("J132").Value = 11183.333
$ws.Range("K132").Value = 14790.5001
$ws.Range("L132").Value = 100649.997
$ws.Range("M132").Value = -12260.5001
$ws.Range("N132").Value = -105709.997
$ws.Range("H139").Value = 1748.6842
$ws.Range("I139").Value = 1326.5625
$ws.Range("K139").Value = 3979.6875
$ws.Range("M139").Value = 1160.3125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2929050
$ws.Range("I10").Value = 6833333.5
$ws.Range("J10").Value = 837.5
$ws.Range("K10").Value = 6833333.5
$ws.Range("L10").Value = 837.5
$ws.Range("M10").Value = -6833164.5
$ws.Range("N10").Value = -1175.5
$ws.Range("H15").Value = 10525.833
$ws.Range("J15").Value = 10525.833
$ws.Range("L15").Value = 10525.833
$ws.Range("N15").Value = -11101.833
$ws.Range("H70").Value = 6591.8
$ws.Range("I70").Value = 6362
$ws.Range("J70").Value = 7051.4
$ws.Range("K70").Value = 6362
$ws.Range("L70").Value = 7051.4
$ws.Range("M70").Value = -6092
$ws.Range("N70").Value = -7591.4
$ws.Range("H73").Value = 6591.8
$ws.Range("I73").Value = 6362
$ws.Range("J73").Value = 7051.4
$ws.Range("K73").Value = 6362
$ws.Range("L73").Value = 7051.4
$ws.Range("M73").Value = -5426
$ws.Range("N73").Value = -8923.4
$ws.Range("H80").Value = 9582.166999999999
$ws.Range("J80").Value = 10799
$ws.Range("L80").Value = 10799
$ws.Range("N80").Value = -12795
$ws.Range("H81").Value = 10525.833
$ws.Range("J81").Value = 10525.833
$ws.Range("L81").Value = 10525.833
$ws.Range("N81").Value = -12521.833
$ws.Range("H83").Value = 9582.166999999999
$ws.Range("J83").Value = 10799
$ws.Range("L83").Value = 53995
$ws.Range("N83").Value = -63979
$ws.Range("H84").Value = 10525.833
$ws.Range("J84").Value = 10525.833
$ws.Range("L84").Value = 31577.499
$ws.Range("N84").Value = -41561.499

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 15159.8
$ws.Range("J11").Value = 19931.666
$ws.Range("L11").Value = 19931.666
$ws.Range("N11").Value = -20215.666
$ws.Range("H74").Value = 23702.75
$ws.Range("J74").Value = 23702.75
$ws.Range("L74").Value = 23702.75
$ws.Range("N74").Value = -25574.75
$ws.Range("H77").Value = 23702.75
$ws.Range("J77").Value = 23702.75
$ws.Range("L77").Value = 71108.25
$ws.Range("N77").Value = -80468.25
